# Generate Report for Handback
#
# The handback status report was regenerated, producing newer timestamps
# for the handoff/handback events and a refreshed translation "Priority"
# value. In the source workbook several cells literally share the same
# text (and therefore the same shared-string entry), so every cell that
# held one of the old values is updated to the corresponding new value.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" cells
# that used to read 2016-09-05 10:18:35.
$wsOverview.Range("G2").Value = "2016-09-05 10:19:39"
$wsOverview.Range("G4").Value = "2016-09-05 10:19:39"
$wsDeDe.Range("H2").Value     = "2016-09-05 10:19:39"
$wsDeDe.Range("H4").Value     = "2016-09-05 10:19:39"

# Priority cells that used to read "ht".
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn Correspond Handoff Datetime (was 2016-09-05 10:18:29).
$wsZhCn.Range("H2").Value = "2016-09-05 10:19:34"
$wsZhCn.Range("H4").Value = "2016-09-05 10:19:34"

# zh-cn Correspond Handback DateTime (was 2016-09-05 10:18:56).
$wsZhCn.Range("K2").Value = "2016-09-05 10:19:53"
$wsZhCn.Range("K4").Value = "2016-09-05 10:19:53"

# de-de Correspond Handback DateTime (was 2016-09-05 10:19:11).
$wsDeDe.Range("K2").Value = "2016-09-05 10:20:01"
$wsDeDe.Range("K4").Value = "2016-09-05 10:20:01"
